# "added payment approval and member addition"
#
# Attendance sheet (sheet1):
#   - Row 2: add Number=9 (A2), change Name from "JJ" to "xc" (B2),
#            add Paid="PAID" (C2) — payment approval for the existing row
#   - Row 3: unchanged values (Number=2, Paid time stamp) — left as-is
#   - Row 4 (new): Number=99 (A4), Name="js" (B4) — new member addition
#   - Row 5: touched so the sheet's used range grows to A1:E5

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Approve payment for the existing attendee in row 2
$ws.Range("A2").Value = 9
$ws.Range("B2").Value = "xc"
$ws.Range("C2").Value = "PAID"

# Add a new member in row 4
$ws.Range("A4").Value = 99
$ws.Range("B4").Value = "js"

# Extend the used range down to row 5 (kept blank)
$ws.Range("A5").Font.Bold = $false
